# Nexial script template update:
#   WEB >> add new command `clickIfPresent(locator)` to the "web" command
#   list on the hidden '#system' sheet (column AE), keeping the list in
#   alphabetical order. The new entry is inserted between
#   `clickByLabelAndWait(label,waitMs)` (row 62) and
#   `clickOffset(locator,x,y)` (row 63), so every existing entry from
#   row 63 down shifts one row lower, and the named range "web" grows
#   from $AE$2:$AE$158 to $AE$2:$AE$159.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$col = 31          # column AE
$insertAt = 63      # new row for "clickIfPresent(locator)"
$lastRow = 158       # last populated row of the "web" list (before the edit)

# Shift existing values in column AE down by one row, starting from the
# bottom so we never overwrite a value before it has been copied down.
# Only column AE is touched -- other columns (e.g. "I", the "desktop"
# list) must stay exactly where they are.
for ($r = $lastRow; $r -ge $insertAt; $r--) {
    $srcCell = $ws.Cells.Item($r, $col)
    $dstCell = $ws.Cells.Item($r + 1, $col)
    $dstCell.Value2 = $srcCell.Value2
}

# Insert the new command in its alphabetically-sorted slot.
$ws.Cells.Item($insertAt, $col).Value2 = "clickIfPresent(locator)"

# Grow the "web" named range so it covers the new row.
$webName = $wb.Names.Item("web")
$webName.RefersTo = '=''#system''!$AE$2:$AE$159'
